$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.461.95"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "3.259.22"
$ws.Range("E3").Value = "  -5.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -11.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.251.30"
$ws.Range("E8").Value = "  -5.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -9.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -13.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.67"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.97%  "
$ws.Range("E12").Value = "  -11.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -8.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -13.94%  "
$ws.Range("D15").Value = "3.776.38"
$ws.Range("E15").Value = "  -5.20%  "
$ws.Range("D16").Value = "67.486.86"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "3.255.57"
$ws.Range("E17").Value = "  -5.07%  "
$ws.Range("E18").Value = "  -6.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "518.66"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -9.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -13.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.73"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -13.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.745"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -11.39%  "
$ws.Range("E23").Value = "  -14.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -10.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.20"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -11.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -12.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.12"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -11.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.62"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -11.88%  "
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -17.15%  "
$ws.Range("E34").Value = "  -14.38%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.31"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "504.75"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -15.45%  "
$ws.Range("E38").Value = "  -6.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0839"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -11.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.80"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -15.59%  "
$ws.Range("E41").Value = "  -11.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -10.53%  "
$ws.Range("D43").Value = "2.913.63"
$ws.Range("E43").Value = "  -9.81%  "
$ws.Range("E44").Value = "  -10.23%  "
$ws.Range("E45").Value = "  -8.90%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.23"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -15.23%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "0.0₃0558"
$ws.Range("E48").Value = "  -16.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.91%  "
$ws.Range("E50").Value = "  -10.50%  "
$ws.Range("E51").Value = "  -18.27%  "
